$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing the existing rows 15-23 down to 16-24
# (mirrors the new weekly price observation added to the dataset).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new observation.
$ws.Cells.Item(15, 1).Value  = 7
$ws.Cells.Item(15, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(15, 3).Value  = 'Ñuble'
$ws.Cells.Item(15, 4).Value  = 45175
$ws.Cells.Item(15, 5).Value  = 16
$ws.Cells.Item(15, 6).Value  = 100112039
$ws.Cells.Item(15, 7).Value  = 'Ciboulette'
$ws.Cells.Item(15, 8).Value  = 'Sin especificar'
$ws.Cells.Item(15, 9).Value  = 'Primera'
$ws.Cells.Item(15, 10).Value = 150
$ws.Cells.Item(15, 11).Value = 2500
$ws.Cells.Item(15, 12).Value = 2500
$ws.Cells.Item(15, 13).Value = 2500
$ws.Cells.Item(15, 14).Value = '$/docena de atados'
$ws.Cells.Item(15, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 16).Value = 833
$ws.Cells.Item(15, 17).Value = 3
$ws.Cells.Item(15, 18).Value = 'Hortaliza'
